# Update gh-pages to output generated at 456a3b4
# Applies numeric "want-to-go" (F) count refresh across the four sheets,
# plus one new exhibition row ("wio夏时之鸢代号鸢Only") inserted into the
# "全部类型" combined sheet, which pushes the existing rows 38-52 down by
# one and drops the oldest trailing row off the bottom of the table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) 展览 (exhibitions) sheet - straight numeric refresh on column F
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1288
$ws.Range("F4").Value = 1330
$ws.Range("F6").Value = 190
$ws.Range("F8").Value = 6
$ws.Range("F9").Value = 42
$ws.Range("F10").Value = 385
$ws.Range("F12").Value = 1306
$ws.Range("F13").Value = 29894
$ws.Range("F14").Value = 5539
$ws.Range("F15").Value = 59
$ws.Range("F16").Value = 282
$ws.Range("F18").Value = 68
$ws.Range("F20").Value = 47
$ws.Range("F21").Value = 33
$ws.Range("F22").Value = 365
$ws.Range("F23").Value = 44
$ws.Range("F24").Value = 677
$ws.Range("F25").Value = 291
$ws.Range("F26").Value = 316
$ws.Range("F27").Value = 376
$ws.Range("F29").Value = 116
$ws.Range("F30").Value = 18
$ws.Range("F31").Value = 686
$ws.Range("F32").Value = 231
$ws.Range("F34").Value = 598
$ws.Range("F35").Value = 83
$ws.Range("F36").Value = 38
$ws.Range("F37").Value = 680
$ws.Range("F38").Value = 254
$ws.Range("F40").Value = 16

# ---------------------------------------------------------------
# 2) 演出 (performances) sheet - straight numeric refresh on column F
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 965
$ws.Range("F7").Value = 1
$ws.Range("F9").Value = 280
$ws.Range("F10").Value = 4280
$ws.Range("F12").Value = 204
$ws.Range("F16").Value = 20
$ws.Range("F22").Value = 4260
$ws.Range("F24").Value = 4

# ---------------------------------------------------------------
# 3) 本地生活 (local-life) sheet - straight numeric refresh on column F
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 304
$ws.Range("F3").Value = 277
$ws.Range("F4").Value = 1277
$ws.Range("F5").Value = 317

# ---------------------------------------------------------------
# 4) 全部类型 (all-types, combined) sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")

# 4a) rows above the insertion point only get the same numeric refresh
$ws.Range("F2").Value = 304
$ws.Range("F3").Value = 277
$ws.Range("F4").Value = 1277
$ws.Range("F7").Value = 317
$ws.Range("F8").Value = 965
$ws.Range("F9").Value = 1288
$ws.Range("F11").Value = 190
$ws.Range("F13").Value = 42
$ws.Range("F14").Value = 385
$ws.Range("F15").Value = 1
$ws.Range("F17").Value = 1306
$ws.Range("F20").Value = 280
$ws.Range("F22").Value = 204
$ws.Range("F23").Value = 204
$ws.Range("F26").Value = 282
$ws.Range("F30").Value = 68
$ws.Range("F31").Value = 47
$ws.Range("F33").Value = 33
$ws.Range("F35").Value = 44
$ws.Range("F36").Value = 677
$ws.Range("F37").Value = 291

# 4b) a new exhibition ("广州·wio夏时之鸢代号鸢Only", already present as row
# 26 of 展览) is now also listed on the combined sheet. Insert a fresh row
# at position 38 - this shifts the previous rows 38-52 down to 39-53 - then
# drop the row that fell off the end of the table (the old row 52, the
# 2024-12-20 小野丽莎 show) so the sheet stays at its original 52 data rows.
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(53).Delete()

# Fill in the newly inserted row 38 with the new event's details.
# (Force text format first so the date-shaped string is not auto-converted
# into a real date serial, matching how the other rows store this column.)
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "2024-08-04"
$ws.Range("C38").Value = "广州·wio夏时之鸢代号鸢Only"
$ws.Range("D38").Value = "黄边三横路一街1号 设计殿堂"
$ws.Range("E38").Value = "2024.08.04 10:00-08.04 17:30"
$ws.Range("F38").Value = 316
$ws.Range("G38").Value = 68.8
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=87434"
$ws.Range("I38").Value = "//i0.hdslb.com/bfs/openplatform/202406/orVoRqXY1718293009879.png"

# Column A is a plain sequential row index (0-based) that is independent of
# which event occupies the row, so restore it for every row the insert
# touched (it gets shifted by the native Insert along with everything else).
For ($r = 38; $r -le 52; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

Write-Output "edit applied"
